$wb = $excel.ActiveWorkbook

# Sheet "展览" (worksheet 1): update "想去人数" (column F) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 366
$wsExpo.Range("F4").Value = 2987
$wsExpo.Range("F6").Value = 624

# Sheet "全部类型" (worksheet 4): update "想去人数" (column F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 366
$wsAll.Range("F6").Value = 2987
$wsAll.Range("F8").Value = 624
